$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of cell -> new text value (values are stored as plain text/inline strings,
# matching the original workbook, not auto-converted numbers/percentages).
$updates = @{
    "D2" = "301.24"
    "E2" = "0.67%"
    "D3" = "31.51"
    "E3" = "0.96%"
    "D4" = "5.083"
    "E4" = "-1.04%"
    "D5" = "0.07843"
    "E5" = "-3.00%"
    "D6" = "2.334"
    "E6" = "-12.75%"
    "D7" = "7.818"
    "E7" = "-0.44%"
    "D8" = "3.834"
    "E8" = "0.02%"
    "D9" = "0.9184"
    "E9" = "1.05%"
    "D10" = "0.1761"
    "E10" = "1.97%"
    "D11" = "0.07590"
    "E11" = "4.49%"
    "D12" = "0.09214"
    "E12" = "14.88%"
    "D13" = "0.02997"
    "E13" = "-0.95%"
    "D14" = "0.1002"
    "E14" = "0.36%"
    "D15" = "0.001521"
    "E15" = "1.55%"
    "D16" = "0.005933"
    "E16" = "-1.76%"
    "D17" = "3.469"
    "E17" = "-0.84%"
    "E18" = "-0.31%"
    "E19" = "-0.51%"
    "D20" = "0.1281"
    "E20" = "-4.39%"
    "D21" = "4.051"
    "E21" = "-12.14%"
    "D22" = "0.1790"
    "E22" = "11.69%"
    "E23" = "0.43%"
    "E24" = "-1.33%"
    "D25" = "0.004471"
    "E25" = "0.60%"
    "E26" = "5.76%"
    "E27" = "-1.55%"
    "D39" = "0.01767"
    "E39" = "-2.49%"
    "D40" = "0.04797"
    "E40" = "5.74%"
    "D41" = "0.007192"
    "E41" = "1.80%"
    "D42" = "0.1359"
    "E42" = "1.12%"
    "D43" = "0.002190"
    "E43" = "-2.39%"
    "D44" = "0.01031"
    "E44" = "-2.06%"
    "D45" = "0.00006351"
    "E45" = "0.27%"
    "E46" = "-0.16%"
    "E47" = "24.60%"
    "D48" = "0.7455"
    "E48" = "-9.15%"
    "D49" = "0.00002100"
    "E49" = "-0.16%"
    "D50" = "0.0002000"
    "E50" = "-0.16%"
}

foreach ($cellRef in $updates.Keys) {
    $value = $updates[$cellRef]
    $range = $ws.Range($cellRef)
    # Leading apostrophe forces text entry (avoids Excel auto-converting
    # numeric-looking strings like "301.24" or "0.67%" into numbers).
    $range.Value = "'" + $value
    # Drop the quote-prefix formatting Excel applies for forced-text entry so
    # the cell keeps the workbook's original (unstyled) appearance.
    $range.ClearFormats()
}
